$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "teste"

$ws.Range("B15").Select()
